$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is a rolling year-over-year series. The oldest year (2009) is
# dropped, every remaining row shifts up by one, and a new row for 2021 is
# appended at the bottom with the same look (bold/centered/bordered year
# label in column A) as the rest of the table.

# 1. Drop the 2009 row (row 2) - this shifts all the other years up by one.
$ws.Rows.Item(2).Delete()

# 2. Copy the formatting of the (now last) year-label cell so the newly
#    appended year label matches the existing style exactly, then fill in
#    the 2021 row of data.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 1060.6005
$ws.Range("C13").Value = 3061.5312
$ws.Range("D13").Value = 445.0918
$ws.Range("E13").Value = 2677.6223
$ws.Range("F13").Value = 1046.7548
$ws.Range("G13").Value = 713.0235
$ws.Range("H13").Value = 176.7345
$ws.Range("I13").Value = 610.7473
$ws.Range("J13").Value = 18925.5264
$ws.Range("K13").Value = 4797.2131
$ws.Range("L13").Value = 29.7287
$ws.Range("M13").Value = 2540.5713
$ws.Range("N13").Value = 6658.6757
$ws.Range("O13").Value = 28380.8889
$ws.Range("P13").Value = 1347.9039
$ws.Range("Q13").Value = 3411.4163
$ws.Range("R13").Value = 5855.3995
$ws.Range("S13").Value = 1989.2449
$ws.Range("T13").Value = 12128.4462
$ws.Range("U13").Value = 4119.3611
